$d = $word.ActiveDocument

$pairs = @(
    @("499×3=", "602×5="),
    @("613×3=", "122×7="),
    @("141×8=", "164×9="),
    @("977×8=", "527×2="),
    @("221×6=", "477×2="),
    @("198×8=", "205×9="),
    @("806×5=", "939×7="),
    @("610×7=", "518×7="),
    @("819×8=", "629×6="),
    @("517×8=", "784×4="),
    @("391×2=", "953×5="),
    @("881×6=", "833×2="),
    @("538×3=", "650×2="),
    @("692×4=", "381×7="),
    @("314×6=", "798×3="),
    @("473×9=", "547×2="),
    @("152×5=", "319×6="),
    @("213×2=", "411×8="),
    @("974×6=", "664×3="),
    @("829×3=", "469×6="),
    @("954×2=", "363×6="),
    @("430×5=", "248×4="),
    @("607×9=", "862×5="),
    @("698×5=", "798×6="),
    @("583×5=", "715×9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
